$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") values for data rows 2-28 are updated from
# serial date 45178 (2023-09-09) to 45179 (2023-09-10).
for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45179
}
